# Duplicate the existing "Shapiro_N12" sheet and place the copy before it.
# This becomes the new first sheet, with the same header row/styles as the
# original, ready to be renamed and re-populated with the N=2,480 data.
$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item("Shapiro_N12")
$orig.Copy($orig)

# The freshly created copy is now the first sheet in the workbook (Excel
# inserts "Copy-Before" targets ahead of the reference sheet).
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Shapiro_N2480"

$ws.Range("A2").Value = "N=2,480"
$ws.Range("B2").Value = "Instruction Coverage (%)"
$ws.Range("C2").Value = "Manual"
$ws.Range("D2").Value = 1600
$ws.Range("E2").Value = [double]"0.816008043429448"
$ws.Range("F2").Value = [double]"2.337136908350902e-39"
$ws.Range("G2").Value = "NO"
$ws.Range("H2").Value = [double]"19.941"
$ws.Range("I2").Value = [double]"12.1149437230795"
$ws.Range("A3").Value = "N=2,480"
$ws.Range("B3").Value = "Instruction Coverage (%)"
$ws.Range("C3").Value = "IA"
$ws.Range("D3").Value = 880
$ws.Range("E3").Value = [double]"0.7523643876549621"
$ws.Range("F3").Value = [double]"2.663606854482367e-34"
$ws.Range("G3").Value = "NO"
$ws.Range("H3").Value = [double]"14.20104545454545"
$ws.Range("I3").Value = [double]"8.707367967776438"
$ws.Range("A4").Value = "N=2,480"
$ws.Range("B4").Value = "Branch Coverage (%)"
$ws.Range("C4").Value = "Manual"
$ws.Range("D4").Value = 1600
$ws.Range("E4").Value = [double]"0.8177699142698152"
$ws.Range("F4").Value = [double]"3.214917981583953e-39"
$ws.Range("G4").Value = "NO"
$ws.Range("H4").Value = [double]"17.6875"
$ws.Range("I4").Value = [double]"12.39259674296598"
$ws.Range("A5").Value = "N=2,480"
$ws.Range("B5").Value = "Branch Coverage (%)"
$ws.Range("C5").Value = "IA"
$ws.Range("D5").Value = 880
$ws.Range("E5").Value = [double]"0.8697179827556469"
$ws.Range("F5").Value = [double]"1.968664763323198e-26"
$ws.Range("G5").Value = "NO"
$ws.Range("H5").Value = [double]"13.50284090909091"
$ws.Range("I5").Value = [double]"6.747887032681644"
$ws.Range("A6").Value = "N=2,480"
$ws.Range("B6").Value = "Mutation Score (%)"
$ws.Range("C6").Value = "Manual"
$ws.Range("D6").Value = 1600
$ws.Range("E6").Value = [double]"0.8404116629778383"
$ws.Range("F6").Value = [double]"2.458931328321721e-37"
$ws.Range("G6").Value = "NO"
$ws.Range("H6").Value = [double]"22.918"
$ws.Range("I6").Value = [double]"17.1925242563866"
$ws.Range("A7").Value = "N=2,480"
$ws.Range("B7").Value = "Mutation Score (%)"
$ws.Range("C7").Value = "IA"
$ws.Range("D7").Value = 880
$ws.Range("E7").Value = [double]"0.8313344102201263"
$ws.Range("F7").Value = [double]"1.820920203199541e-29"
$ws.Range("G7").Value = "NO"
$ws.Range("H7").Value = [double]"16.63388636363636"
$ws.Range("I7").Value = [double]"8.160557021954679"
$ws.Range("A8").Value = "N=2,480"
$ws.Range("B8").Value = "Time (seconds)"
$ws.Range("C8").Value = "Manual"
$ws.Range("D8").Value = 1600
$ws.Range("E8").Value = [double]"0.479676697019995"
$ws.Range("F8").Value = [double]"8.07180512008333e-56"
$ws.Range("G8").Value = "NO"
$ws.Range("H8").Value = [double]"0.079395625"
$ws.Range("I8").Value = [double]"0.1775453021479479"
$ws.Range("A9").Value = "N=2,480"
$ws.Range("B9").Value = "Time (seconds)"
$ws.Range("C9").Value = "IA"
$ws.Range("D9").Value = 880
$ws.Range("E9").Value = [double]"0.5354067631882183"
$ws.Range("F9").Value = [double]"4.671058967064643e-43"
$ws.Range("G9").Value = "NO"
$ws.Range("H9").Value = [double]"0.1140261363636364"
$ws.Range("I9").Value = [double]"0.231760935630358"

# Make the new sheet the active tab, matching activeTab="0" in the workbook.
$ws.Activate()
